$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "57.395.69"
$ws.Range("E2").Value = "  +1.66%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.072.64"
$ws.Range("E3").Value = "  +2.84%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "515.06"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.06%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.67"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.84%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.14%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +2.06%  "

# Row 9 - Toncoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.32"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.18%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.63%  "

# Row 11 - Cardano
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.374"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.45%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "3.589.49"
$ws.Range("E12").Value = "  +4.04%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +2.91%  "

# Row 14 - Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.66"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.09%  "

# Row 15 - ShibaInu
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000165"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.16%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "57.499.20"
$ws.Range("E16").Value = "  +1.75%  "

# Row 17 - WrappedEther
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.13"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.02%  "

# Row 18 - Polkadot
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.065.22"
$ws.Range("E18").Value = "  +3.90%  "

# Row 19 - Chainlink
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.00"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.11%  "

# Row 20 - Uniswap
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.18"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.21%  "

# Row 21 - BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "337.11"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.76%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.24%  "

# Row 23 - Polygon
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.500"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.00%  "

# Row 24 - Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.41"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.11%  "

# Row 25 - Kaspa
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.172"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +6.13%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("E26").Value = "  -0.03%  "

# Row 27 - PEPE
$ws.Range("D27").Value = "0.0₃0930"
$ws.Range("E27").Value = "  +6.86%  "

# Row 28 - RenderToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.43"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.93%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.06"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.26%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +1.21%  "

# Row 31 - EthereumClassic
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.83"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.46%  "

# Row 32 - Fetch.AI
$ws.Range("E32").Value = "  -0.34%  "

# Row 33 - Monero
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.13"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.17%  "

# Row 34 - NEARProtocol
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.52"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.08%  "

# Row 35 - Aptos
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.90"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.47%  "

# Row 36 - EnergySwap
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.55"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +7.42%  "

# Row 37 - ImmutableX
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.24"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.94%  "

# Row 38 - Hedera
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0684"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.85%  "

# Row 39 - RenzoRestakedETH
$ws.Range("D39").Value = "3.106.47"
$ws.Range("E39").Value = "  +3.72%  "

# Row 40 - OKB
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.89"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.89%  "

# Row 41 - Filecoin
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.87"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.98%  "

# Row 42 - FirstDigitalUSD
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.669"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.93%  "

# Row 43 - Mantle
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.10%  "

# Row 44 - Maker
$ws.Range("D44").Value = "2.263.69"
$ws.Range("E44").Value = "  +5.49%  "

# Row 45 - VeChain
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0253"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +7.23%  "

# Row 46 - Stacks
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.39"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.84%  "

# Row 47 - ONDO
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.957"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.09%  "

# Row 48 - InjectiveProtocol
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.34"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +6.35%  "

# Row 49 - Cosmos
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.86"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.12%  "

# Row 50 - Stellar
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0875"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.67%  "

# Row 51 - SuiNetwork
$ws.Range("E51").Value = "  +4.89%  "
